$d = $word.ActiveDocument

# Locate the paragraph that ends the "GET /static/Logo.jpg" log line -- the
# new block of paragraphs must be inserted right after it (and right before
# the paragraph that currently carries the bottom border separator).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains('GET /static/Logo.jpg HTTP/1.1" 200 18645')) {
        $anchor = $p
    }
}

if ($anchor -eq $null) {
    throw "Could not find the anchor paragraph (GET /static/Logo.jpg line)."
}

# Zero-length insertion point right after the anchor paragraph's own mark,
# i.e. exactly between it and the following (border) paragraph.
$pos = $anchor.Range.End
$ins = $d.Range($pos, $pos)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$fragment = "<w:p $ns>" +
                "<w:pPr>" +
                    "<w:pBdr><w:bottom w:val='single' w:sz='6' w:space='1' w:color='auto'/></w:pBdr>" +
                    "<w:spacing w:after='0'/>" +
                "</w:pPr>" +
            "</w:p>" +
            "<w:p $ns>" +
                "<w:pPr><w:spacing w:after='0'/></w:pPr>" +
                "<w:r><w:tab/></w:r>" +
                "<w:r><w:tab/></w:r>" +
                "<w:r><w:tab/></w:r>" +
                "<w:r><w:tab/><w:t>** Dynamic Data and Pic. Added 9**</w:t></w:r>" +
            "</w:p>" +
            "<w:p $ns>" +
                "<w:pPr><w:spacing w:after='0'/></w:pPr>" +
            "</w:p>" +
            "<w:p $ns>" +
                "<w:pPr><w:spacing w:after='0'/></w:pPr>" +
            "</w:p>"

$null = $ins.InsertXML($fragment)
